$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "karthika"
$ws.Range("F2").Value = "vijayan"

$ws.Range("F2").Select()
